$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Range("E11").Value = 0.3629913849272427
$ws.Range("E12").Value = 0.12488235719973945
$ws.Range("G12").Value = 3162.5000000000005
$ws.Range("E13").Value = 0.30866755954535602
$ws.Range("E14").Value = 0.18732353579960917
$ws.Range("G14").Value = 3162.5000000000009
$ws.Range("E150").Value = 0.16085025627375071
$ws.Range("E151").Value = 0.15226887751132734
$ws.Range("E152").Value = 0.15456128021356608
$ws.Range("E153").Value = 0.16629376698088194
$ws.Range("E154").Value = 0.17206733071733712
$ws.Range("E155").Value = 0.15273795001145538
$ws.Range("E156").Value = 0.21640319337561012
$ws.Range("E157").Value = 0.19699531868281184
$ws.Range("E158").Value = 0.19247860444770779
$ws.Range("E159").Value = 0.19228757088918788
$ws.Range("E160").Value = 0.21063530390326943
$ws.Range("E161").Value = 0.13549669849969209
$ws.Range("E162").Value = 0.18231505170803797
$ws.Range("E163").Value = 0.13652468601509371
$ws.Range("E164").Value = 0.20006982412215921
$ws.Range("E165").Value = 0.21381383751804844
$ws.Range("E166").Value = 0.16209575724687297
$ws.Range("E167").Value = 0.1579516530219513
$ws.Range("E168").Value = 0.16568094645652107
$ws.Range("E169").Value = 0.1534481787364477
$ws.Range("E170").Value = 0.19745398836539674
$ws.Range("E171").Value = 0.20988535532947597
$ws.Range("E172").Value = 0.1701049810444224
$ws.Range("E173").Value = 0.19532613932247714
$ws.Range("E174").Value = 0.16485344960649678
